# Work-log.xlsx update
# Adds two new log rows (39 and 40) for "python/pip" to Sheet1, mirroring
# the formatting of the preceding row (38), adjusts column widths for
# columns B and C, and updates the selection to the new last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (fill/border/number format/font) of the last existing
# data row (38) down onto the two new rows so the new cells pick up the same
# cell styles already present in the workbook (date format on col A, etc.)
# instead of creating brand-new style entries.
$ws.Range("A38:C38").Copy()
$ws.Range("A39:C40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode() = $false

# Row 39: python/pip - pip intro & basics
$ws.Range("A39").Value() = 43348
$ws.Range("B39").Value() = "python/pip"
$ws.Range("C39").Value() = "pip-intro,installation,requirement files,contraint files,basic commands of pip"

# Row 40: python/pip - packaging task
$ws.Range("A40").Value() = 43348
$ws.Range("B40").Value() = "python/pip"
$ws.Range("C40").Value() = "Created Package and Published the package and then installed it using pip"

# Widen/narrow columns B and C to fit the new content.
$ws.Columns.Item(2).ColumnWidth() = 30.166666666666668
$ws.Columns.Item(3).ColumnWidth() = 118.30729166666667

# Move the selection to the new bottom-right cell, like the author did.
$ws.Range("C40").Select()
